$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated metric values for B2:D10 (MSE, R2, MAE) ---
$values = @(
    @(1.332598112885935, 0.7829060121382954, 0.9583714877187218),
    @(6.90266466986186,  0.901389913661654,  1.983252460510221),
    @(3.377229054509587, 0.8331785694834976, 1.481931750261225),
    @(3.095584561444309, 0.9979666919584735, 1.291599192029326),
    @(2.127566867576145, 0.9768165028372422, 1.162564383981743),
    @(1.817672486595947, 0.9989808878376261, 1.065162419437801),
    @(2.273250848686807, 0.9973960973649652, 1.248790469195605),
    @(15.49277267495131, 0.8147946083240025, 3.264842200283293),
    @(1.535623997710418, 0.9954223524263063, 1.001821233815439)
)

# --- New "Modelo" column header (F1) ---
$ws.Range("F1").Value = "Modelo"
# Copy format from an existing header cell so F1 gets the same bold/border/center style
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Model description text (multi-line) written via a formula using CHAR(10) ---
# then converted to a plain value, so that Excel does not auto-resize the row
# height to fit the embedded newlines (matches the unedited row heights).
$lines = @(
    "MultiOutputRegressor(estimator=GridSearchCV(cv=5,",
    "                                            estimator=Pipeline(steps=[('model',",
    "                                                                       RandomForestRegressor())]),",
    "                                            param_grid={'model__max_depth': [3,",
    "                                                                             5,",
    "                                                                             7],",
    "                                                        'model__n_estimators': [50,",
    "                                                                                100,",
    "                                                                                150]},",
    "                                            scoring='neg_mean_squared_error'))"
)

$parts = @()
foreach ($l in $lines) {
    $escaped = $l.Replace('"', '""')
    $parts += ('"' + $escaped + '"')
}
$modelFormula = "=" + [string]::Join("&CHAR(10)&", $parts)

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
    $ws.Cells.Item($row, 4).Value = $values[$i][2]

    $cell = $ws.Cells.Item($row, 6)
    $cell.Formula = $modelFormula
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues: bake formula result into a static value
}
$excel.CutCopyMode = 0
